# BreakoutBoard_BOM.xlsx edit script
#
# Summary of the change (per the commit's xml diff):
#   - Drop "Sheet2" and "Sheet3" (now-unused blank sheets).
#   - Rename the remaining sheet from "Power Supply Board - Rev A" to
#     "Breakout Board - Rev B".
#   - Swap the Item # values for the two header rows (J1/J2 connectors).
#   - Update the title banner text in A1 to the new BOM title.
#   - Change the active selection from B17 to A1:L1.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Power Supply Board - Rev A")

# Remove the two blank, unused worksheets.
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Sheet2").Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheet to match the new board revision.
$ws.Name = "Breakout Board - Rev B"

# Swap the Item # numbering between the J1 and J2 connector rows.
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2

# Update the BOM title text.
$ws.Range("A1").Value = "Bill of Materials for 'Marmote - Breakout Board Rev B (Smoky)'"

# Update the selection/active range on the sheet.
[void]$ws.Range("A1:L1").Select()
